$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 3 cells per diff (odds changed) ---
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 11
$ws.Range("Q3").Value = 2
$ws.Range("R3").Value = 1.85

# --- New row 4 ---
$ws.Range("A4").Value = "dximVNgc"
$ws.Range("B4").Value = "25/10/2024"
$ws.Range("C4").Value = "09:15"
$ws.Range("D4").Value = "MALAYSIA - SUPER LEAGUE"
$ws.Range("E4").Value = "Penang"
$ws.Range("F4").Value = "Sabah"
$ws.Range("G4").Value = 2.67
$ws.Range("H4").Value = 3.5
$ws.Range("I4").Value = 2.22
$ws.Range("J4").Value = 3.15
$ws.Range("K4").Value = 2.25
$ws.Range("L4").Value = 2.7
$ws.Range("M4").Value = 1.03
$ws.Range("N4").Value = 10
$ws.Range("O4").Value = 1.17
$ws.Range("P4").Value = 4.34
$ws.Range("Q4").Value = 1.6
$ws.Range("R4").Value = 2.07
$ws.Range("S4").Value = 1.3
$ws.Range("T4").Value = 3.3
$ws.Range("U4").Value = 1.52
$ws.Range("V4").Value = 2.46
$ws.Range("W4").Value = 9.5
$ws.Range("X4").Value = 13
$ws.Range("Y4").Value = 8.5
$ws.Range("Z4").Value = 25
$ws.Range("AA4").Value = 16.5
$ws.Range("AB4").Value = 20
$ws.Range("AC4").Value = 13.5
$ws.Range("AD4").Value = 6.3
$ws.Range("AE4").Value = 10.25
$ws.Range("AF4").Value = 32
$ws.Range("AG4").Value = 175
$ws.Range("AH4").Value = 8.75
$ws.Range("AI4").Value = 10.75
$ws.Range("AJ4").Value = 7.7
$ws.Range("AK4").Value = 18.5
$ws.Range("AL4").Value = 13.5
$ws.Range("AM4").Value = 17.5
$ws.Range("AN4").Value = 4.9
$ws.Range("AO4").Value = 14
$ws.Range("AP4").Value = 18.5
$ws.Range("AQ4").Value = 55
$ws.Range("AR4").Value = 75
$ws.Range("AS4").Value = 200
$ws.Range("AT4").Value = 3.2
$ws.Range("AU4").Value = 6.5
$ws.Range("AV4").Value = 45
$ws.Range("AW4").Value = 51
$ws.Range("AX4").Value = 4.4
$ws.Range("AY4").Value = 11
$ws.Range("AZ4").Value = 16.5
$ws.Range("BA4").Value = 40
$ws.Range("BB4").Value = 60
$ws.Range("BC4").Value = 175
$ws.Range("BD4").Value = 51

# --- New row 5 ---
$ws.Range("A5").Value = "Y1UhtUYP"
$ws.Range("B5").Value = "25/10/2024"
$ws.Range("C5").Value = "09:30"
$ws.Range("D5").Value = "UKRAINE - PREMIER LEAGUE"
$ws.Range("E5").Value = "Rukh Lviv"
$ws.Range("F5").Value = "Ch. Odesa"
$ws.Range("G5").Value = 1.72
$ws.Range("H5").Value = 3.25
$ws.Range("I5").Value = 4.9
$ws.Range("J5").Value = 2.37
$ws.Range("K5").Value = 1.95
$ws.Range("L5").Value = 5.4
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7
$ws.Range("O5").Value = 1.47
$ws.Range("P5").Value = 2.32
$ws.Range("Q5").Value = 2.37
$ws.Range("R5").Value = 1.45
$ws.Range("S5").Value = 1.52
$ws.Range("T5").Value = 2.22
$ws.Range("U5").Value = 2.22
$ws.Range("V5").Value = 1.52
$ws.Range("W5").Value = 4.9
$ws.Range("X5").Value = 6.7
$ws.Range("Y5").Value = 9
$ws.Range("Z5").Value = 13
$ws.Range("AA5").Value = 17.5
$ws.Range("AB5").Value = 45
$ws.Range("AC5").Value = 6.7
$ws.Range("AD5").Value = 6.7
$ws.Range("AE5").Value = 23
$ws.Range("AF5").Value = 175
$ws.Range("AG5").Value = 67
$ws.Range("AH5").Value = 9.75
$ws.Range("AI5").Value = 26
$ws.Range("AJ5").Value = 18
$ws.Range("AK5").Value = 100
$ws.Range("AL5").Value = 70
$ws.Range("AM5").Value = 90
$ws.Range("AN5").Value = 3.3
$ws.Range("AO5").Value = 8.75
$ws.Range("AP5").Value = 24
$ws.Range("AQ5").Value = 32
$ws.Range("AR5").Value = 90
$ws.Range("AS5").Value = 450
$ws.Range("AT5").Value = 2.2
$ws.Range("AU5").Value = 8.75
$ws.Range("AV5").Value = 120
$ws.Range("AW5").Value = 81
$ws.Range("AX5").Value = 6.3
$ws.Range("AY5").Value = 32
$ws.Range("AZ5").Value = 45
$ws.Range("BA5").Value = 250
$ws.Range("BB5").Value = 350
$ws.Range("BC5").Value = 700
$ws.Range("BD5").Value = 81

